$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = 47.25208013818724
$ws.Range("D12").Value = 25.53591645352849
$ws.Range("D23").Value = 43.53639735871651
$ws.Range("D34").Value = 11.0700097266827
$ws.Range("E34").Value = 18
